$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("username") to hold the new "browser" data
$ws.Columns.Item(3).Insert()

# Populate the new "browser" column (C)
$ws.Range("C1").Value = "browser"
$ws.Range("C2").Value = "chrome"
$ws.Range("C3").Value = "firefox"
$ws.Range("C4").Value = "'"
$ws.Range("C5").Value = "'"

# Update execute flags in column B to reflect the new cross-browser rows
$ws.Range("B3").Value = "yes"
$ws.Range("B4").Value = "no"
$ws.Range("B5").Value = "no"

# Move the active selection to C6
$ws.Range("C6").Select()
